$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D182").Value = 44476
$ws.Range("L182").Value = "1a nueva(o)"
$ws.Range("M182").Value = 150
$ws.Range("N182").Value = 4000
$ws.Range("O182").Value = 4000
$ws.Range("P182").Value = 4000
$ws.Range("Q182").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R182").Value = "Provincia de Quillota"
$ws.Range("S182").Value = 4000
$ws.Range("T182").Value = 1

$ws.Range("D183").Value = 44476
$ws.Range("L183").Value = "2a nueva(o)"
$ws.Range("M183").Value = 150
$ws.Range("N183").Value = 3500
$ws.Range("O183").Value = 3500
$ws.Range("P183").Value = 3500
$ws.Range("Q183").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R183").Value = "Provincia de Quillota"
$ws.Range("S183").Value = 3500
$ws.Range("T183").Value = 1

$ws.Range("D184").Value = 44386
$ws.Range("L184").Value = "1a nueva(o)"
$ws.Range("M184").Value = 300
$ws.Range("N184").Value = 4500
$ws.Range("O184").Value = 4800
$ws.Range("P184").Value = 4650
$ws.Range("Q184").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R184").Value = "Región de O'Higgins"
$ws.Range("S184").Value = 4650
$ws.Range("T184").Value = 1

$ws.Range("D185").Value = 44386
$ws.Range("L185").Value = "2a nueva(o)"
$ws.Range("M185").Value = 150
$ws.Range("N185").Value = 4000
$ws.Range("O185").Value = 4000
$ws.Range("P185").Value = 4000
$ws.Range("Q185").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R185").Value = "Región de O'Higgins"
$ws.Range("S185").Value = 4000
$ws.Range("T185").Value = 1

$ws.Range("D186").Value = 44306
$ws.Range("L186").Value = "Primera"
$ws.Range("M186").Value = 200
$ws.Range("N186").Value = 45000
$ws.Range("O186").Value = 46000
$ws.Range("P186").Value = 45500
$ws.Range("Q186").Value = "`$/bandeja 10 kilos"
$ws.Range("R186").Value = "Perú"
$ws.Range("S186").Value = 4550
$ws.Range("T186").Value = 10

$ws.Range("D187").Value = 44306
$ws.Range("L187").Value = "Primera"
$ws.Range("M187").Value = 200
$ws.Range("N187").Value = 6900
$ws.Range("O187").Value = 7000
$ws.Range("P187").Value = 6950
$ws.Range("Q187").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R187").Value = "Provincia de Quillota"
$ws.Range("S187").Value = 6950
$ws.Range("T187").Value = 1

$ws.Range("D188").Value = 44369
$ws.Range("L188").Value = "1a nueva(o)"
$ws.Range("M188").Value = 300
$ws.Range("N188").Value = 5000
$ws.Range("O188").Value = 5200
$ws.Range("P188").Value = 5100
$ws.Range("Q188").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R188").Value = "Provincia de Quillota"
$ws.Range("S188").Value = 5100
$ws.Range("T188").Value = 1

$ws.Range("D189").Value = 44369
$ws.Range("L189").Value = "2a nueva(o)"
$ws.Range("M189").Value = 150
$ws.Range("N189").Value = 4200
$ws.Range("O189").Value = 4200
$ws.Range("P189").Value = 4200
$ws.Range("Q189").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R189").Value = "Provincia de Quillota"
$ws.Range("S189").Value = 4200
$ws.Range("T189").Value = 1

$ws.Range("D190").Value = 44369
$ws.Range("L190").Value = "Primera"
$ws.Range("M190").Value = 150
$ws.Range("N190").Value = 40000
$ws.Range("O190").Value = 40000
$ws.Range("P190").Value = 40000
$ws.Range("Q190").Value = "`$/bandeja 10 kilos"
$ws.Range("R190").Value = "Perú"
$ws.Range("S190").Value = 4000
$ws.Range("T190").Value = 10

$ws.Range("D191").Value = 44172
$ws.Range("L191").Value = "Primera"
$ws.Range("M191").Value = 100
$ws.Range("N191").Value = 4500
$ws.Range("O191").Value = 4500
$ws.Range("P191").Value = 4500
$ws.Range("Q191").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R191").Value = "Provincia de Quillota"
$ws.Range("S191").Value = 4500
$ws.Range("T191").Value = 1

$ws.Range("D192").Value = 44172
$ws.Range("L192").Value = "Segunda"
$ws.Range("M192").Value = 100
$ws.Range("N192").Value = 3800
$ws.Range("O192").Value = 3800
$ws.Range("P192").Value = 3800
$ws.Range("Q192").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R192").Value = "Provincia de Quillota"
$ws.Range("S192").Value = 3800
$ws.Range("T192").Value = 1

$ws.Range("D193").Value = 44172
$ws.Range("L193").Value = "Tercera"
$ws.Range("M193").Value = 100
$ws.Range("N193").Value = 3000
$ws.Range("O193").Value = 3000
$ws.Range("P193").Value = 3000
$ws.Range("Q193").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R193").Value = "Provincia de Quillota"
$ws.Range("S193").Value = 3000
$ws.Range("T193").Value = 1

$ws.Range("D194").Value = 44223
$ws.Range("L194").Value = "Primera"
$ws.Range("M194").Value = 100
$ws.Range("N194").Value = 5600
$ws.Range("O194").Value = 5700
$ws.Range("P194").Value = 5650
$ws.Range("Q194").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R194").Value = "Provincia de Quillota"
$ws.Range("S194").Value = 5650
$ws.Range("T194").Value = 1

$ws.Range("D195").Value = 44223
$ws.Range("L195").Value = "Segunda"
$ws.Range("M195").Value = 60
$ws.Range("N195").Value = 4700
$ws.Range("O195").Value = 4700
$ws.Range("P195").Value = 4700
$ws.Range("Q195").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R195").Value = "Provincia de Quillota"
$ws.Range("S195").Value = 4700
$ws.Range("T195").Value = 1

$ws.Range("D196").Value = 44298
$ws.Range("L196").Value = "Primera"
$ws.Range("M196").Value = 80
$ws.Range("N196").Value = 6400
$ws.Range("O196").Value = 6500
$ws.Range("P196").Value = 6450
$ws.Range("Q196").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R196").Value = "Provincia de Quillota"
$ws.Range("S196").Value = 6450
$ws.Range("T196").Value = 1

$ws.Range("D197").Value = 44397
$ws.Range("L197").Value = "Especial"
$ws.Range("M197").Value = 200
$ws.Range("N197").Value = 40000
$ws.Range("O197").Value = 40000
$ws.Range("P197").Value = 40000
$ws.Range("Q197").Value = "`$/bandeja 10 kilos"
$ws.Range("R197").Value = "Perú"
$ws.Range("S197").Value = 4000
$ws.Range("T197").Value = 10

$ws.Range("D198").Value = 44414
$ws.Range("L198").Value = "Especial"
$ws.Range("M198").Value = 150
$ws.Range("N198").Value = 35000
$ws.Range("O198").Value = 35000
$ws.Range("P198").Value = 35000
$ws.Range("Q198").Value = "`$/bandeja 10 kilos"
$ws.Range("R198").Value = "Perú"
$ws.Range("S198").Value = 3500
$ws.Range("T198").Value = 10

$ws.Range("D199").Value = 44414
$ws.Range("L199").Value = "Primera"
$ws.Range("M199").Value = 300
$ws.Range("N199").Value = 28000
$ws.Range("O199").Value = 28000
$ws.Range("P199").Value = 28000
$ws.Range("Q199").Value = "`$/bandeja 10 kilos"
$ws.Range("R199").Value = "Perú"
$ws.Range("S199").Value = 2800
$ws.Range("T199").Value = 10

$ws.Range("D200").Value = 44314
$ws.Range("L200").Value = "Primera"
$ws.Range("M200").Value = 80
$ws.Range("N200").Value = 7000
$ws.Range("O200").Value = 7200
$ws.Range("P200").Value = 7100
$ws.Range("Q200").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R200").Value = "Provincia de Quillota"
$ws.Range("S200").Value = 7100
$ws.Range("T200").Value = 1

$ws.Range("D201").Value = 44392
$ws.Range("L201").Value = "2a nueva(o)"
$ws.Range("M201").Value = 300
$ws.Range("N201").Value = 3800
$ws.Range("O201").Value = 3800
$ws.Range("P201").Value = 3800
$ws.Range("Q201").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R201").Value = "Provincia de Quillota"
$ws.Range("S201").Value = 3800
$ws.Range("T201").Value = 1

$ws.Range("D202").Value = 44425
$ws.Range("L202").Value = "Especial"
$ws.Range("M202").Value = 100
$ws.Range("N202").Value = 35000
$ws.Range("O202").Value = 35000
$ws.Range("P202").Value = 35000
$ws.Range("Q202").Value = "`$/bandeja 10 kilos"
$ws.Range("R202").Value = "Perú"
$ws.Range("S202").Value = 3500
$ws.Range("T202").Value = 10

$ws.Range("D203").Value = 44425
$ws.Range("L203").Value = "Primera"
$ws.Range("M203").Value = 100
$ws.Range("N203").Value = 32000
$ws.Range("O203").Value = 32000
$ws.Range("P203").Value = 32000
$ws.Range("Q203").Value = "`$/bandeja 10 kilos"
$ws.Range("R203").Value = "Perú"
$ws.Range("S203").Value = 3200
$ws.Range("T203").Value = 10

$ws.Range("D204").Value = 44425
$ws.Range("L204").Value = "Segunda"
$ws.Range("M204").Value = 100
$ws.Range("N204").Value = 24000
$ws.Range("O204").Value = 24000
$ws.Range("P204").Value = 24000
$ws.Range("Q204").Value = "`$/bandeja 10 kilos"
$ws.Range("R204").Value = "Perú"
$ws.Range("S204").Value = 2400
$ws.Range("T204").Value = 10

$ws.Range("D205").Value = 44390
$ws.Range("L205").Value = "2a nueva(o)"
$ws.Range("M205").Value = 350
$ws.Range("N205").Value = 4000
$ws.Range("O205").Value = 4000
$ws.Range("P205").Value = 4000
$ws.Range("Q205").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R205").Value = "Provincia de Quillota"
$ws.Range("S205").Value = 4000
$ws.Range("T205").Value = 1

$ws.Range("D206").Value = 44187
$ws.Range("L206").Value = "Primera"
$ws.Range("M206").Value = 150
$ws.Range("N206").Value = 5100
$ws.Range("O206").Value = 5100
$ws.Range("P206").Value = 5100
$ws.Range("Q206").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R206").Value = "Provincia de Quillota"
$ws.Range("S206").Value = 5100
$ws.Range("T206").Value = 1

$ws.Range("D207").Value = 44187
$ws.Range("L207").Value = "Segunda"
$ws.Range("M207").Value = 150
$ws.Range("N207").Value = 4300
$ws.Range("O207").Value = 4300
$ws.Range("P207").Value = 4300
$ws.Range("Q207").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R207").Value = "Provincia de Quillota"
$ws.Range("S207").Value = 4300
$ws.Range("T207").Value = 1

$ws.Range("D208").Value = 44187
$ws.Range("L208").Value = "Tercera"
$ws.Range("M208").Value = 150
$ws.Range("N208").Value = 3600
$ws.Range("O208").Value = 3600
$ws.Range("P208").Value = 3600
$ws.Range("Q208").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R208").Value = "Provincia de Quillota"
$ws.Range("S208").Value = 3600
$ws.Range("T208").Value = 1

$ws.Range("D209").Value = 44466
$ws.Range("L209").Value = "1a nueva(o)"
$ws.Range("M209").Value = 100
$ws.Range("N209").Value = 3500
$ws.Range("O209").Value = 3500
$ws.Range("P209").Value = 3500
$ws.Range("Q209").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R209").Value = "Provincia de Quillota"
$ws.Range("S209").Value = 3500
$ws.Range("T209").Value = 1

$ws.Range("D210").Value = 44250
$ws.Range("L210").Value = "Primera"
$ws.Range("M210").Value = 200
$ws.Range("N210").Value = 5400
$ws.Range("O210").Value = 5500
$ws.Range("P210").Value = 5450
$ws.Range("Q210").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R210").Value = "Provincia de Quillota"
$ws.Range("S210").Value = 5450
$ws.Range("T210").Value = 1

$ws.Range("D211").Value = 44250
$ws.Range("L211").Value = "Segunda"
$ws.Range("M211").Value = 100
$ws.Range("N211").Value = 4300
$ws.Range("O211").Value = 4300
$ws.Range("P211").Value = 4300
$ws.Range("Q211").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R211").Value = "Provincia de Quillota"
$ws.Range("S211").Value = 4300
$ws.Range("T211").Value = 1

$ws.Range("D212").Value = 44201
$ws.Range("L212").Value = "Primera"
$ws.Range("M212").Value = 300
$ws.Range("N212").Value = 5400
$ws.Range("O212").Value = 5600
$ws.Range("P212").Value = 5500
$ws.Range("Q212").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R212").Value = "Provincia de Quillota"
$ws.Range("S212").Value = 5500
$ws.Range("T212").Value = 1

$ws.Range("D213").Value = 44201
$ws.Range("L213").Value = "Segunda"
$ws.Range("M213").Value = 150
$ws.Range("N213").Value = 4500
$ws.Range("O213").Value = 4500
$ws.Range("P213").Value = 4500
$ws.Range("Q213").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R213").Value = "Provincia de Quillota"
$ws.Range("S213").Value = 4500
$ws.Range("T213").Value = 1

$ws.Range("D214").Value = 44193
$ws.Range("L214").Value = "Primera"
$ws.Range("M214").Value = 200
$ws.Range("N214").Value = 5300
$ws.Range("O214").Value = 5400
$ws.Range("P214").Value = 5350
$ws.Range("Q214").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R214").Value = "Provincia de Quillota"
$ws.Range("S214").Value = 5350
$ws.Range("T214").Value = 1

$ws.Range("D215").Value = 44193
$ws.Range("L215").Value = "Segunda"
$ws.Range("M215").Value = 100
$ws.Range("N215").Value = 4500
$ws.Range("O215").Value = 4500
$ws.Range("P215").Value = 4500
$ws.Range("Q215").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R215").Value = "Provincia de Quillota"
$ws.Range("S215").Value = 4500
$ws.Range("T215").Value = 1

$ws.Range("D216").Value = 44286
$ws.Range("L216").Value = "Primera"
$ws.Range("M216").Value = 80
$ws.Range("N216").Value = 6300
$ws.Range("O216").Value = 6400
$ws.Range("P216").Value = 6350
$ws.Range("Q216").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R216").Value = "Provincia de Quillota"
$ws.Range("S216").Value = 6350
$ws.Range("T216").Value = 1

$ws.Range("D217").Value = 44326
$ws.Range("L217").Value = "Primera"
$ws.Range("M217").Value = 80
$ws.Range("N217").Value = 7000
$ws.Range("O217").Value = 7200
$ws.Range("P217").Value = 7100
$ws.Range("Q217").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R217").Value = "Provincia de Quillota"
$ws.Range("S217").Value = 7100
$ws.Range("T217").Value = 1

$ws.Range("A218").Value = 4
$ws.Range("B218").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C218").Value = "Los Lagos"
$ws.Range("D218").Value = 44432
$ws.Range("D218").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E218").Value = 10
$ws.Range("F218").Value = "Fruta"
$ws.Range("G218").Value = 100106
$ws.Range("H218").Value = "Oleaginosos"
$ws.Range("I218").Value = 100106002
$ws.Range("J218").Value = "Palta"
$ws.Range("K218").Value = "Hass"
$ws.Range("L218").Value = "Especial"
$ws.Range("M218").Value = 150
$ws.Range("N218").Value = 35000
$ws.Range("O218").Value = 35000
$ws.Range("P218").Value = 35000
$ws.Range("Q218").Value = "`$/bandeja 10 kilos"
$ws.Range("R218").Value = "Perú"
$ws.Range("S218").Value = 3500
$ws.Range("T218").Value = 10

$ws.Range("A219").Value = 4
$ws.Range("B219").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C219").Value = "Los Lagos"
$ws.Range("D219").Value = 44432
$ws.Range("D219").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E219").Value = 10
$ws.Range("F219").Value = "Fruta"
$ws.Range("G219").Value = 100106
$ws.Range("H219").Value = "Oleaginosos"
$ws.Range("I219").Value = 100106002
$ws.Range("J219").Value = "Palta"
$ws.Range("K219").Value = "Hass"
$ws.Range("L219").Value = "Primera"
$ws.Range("M219").Value = 150
$ws.Range("N219").Value = 30000
$ws.Range("O219").Value = 30000
$ws.Range("P219").Value = 30000
$ws.Range("Q219").Value = "`$/bandeja 10 kilos"
$ws.Range("R219").Value = "Perú"
$ws.Range("S219").Value = 3000
$ws.Range("T219").Value = 10
